$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new blank column before column N (14) ---
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = 9.83

# Move the active selection / active sheet to match the edited workbook state
$ws.Range("P7").Select() | Out-Null
$wb.Worksheets.Item("Transactions").Select() | Out-Null
$ws.Select() | Out-Null
